# Add New Client feature
# Inserts a "Phone Number" column ahead of "Customer ID", backfills the
# existing rows with blank placeholders + generated Customer IDs, and
# appends a new client row (Aaron Gallaway) at the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Client List")

# --- Insert the new "Phone Number" column before the existing "Customer ID" column (H) ---
$ws.Columns.Item(8).Insert()
$ws.Cells.Item(1,8).Value = "Phone Number"

# --- Backfill existing client rows (2-7): blank Address/City/State/Zip/Phone
#     placeholders plus a generated numeric Customer ID in the new last column (I) ---
$ws.Cells.Item(2,4).Value = " "
$ws.Cells.Item(2,5).Value = " "
$ws.Cells.Item(2,6).Value = " "
$ws.Cells.Item(2,7).Value = " "
$ws.Cells.Item(2,8).Value = " "
$ws.Cells.Item(2,9).Value = 814535

$ws.Cells.Item(3,4).Value = " "
$ws.Cells.Item(3,5).Value = " "
$ws.Cells.Item(3,6).Value = " "
$ws.Cells.Item(3,7).Value = " "
$ws.Cells.Item(3,8).Value = " "
$ws.Cells.Item(3,9).Value = 733079

$ws.Cells.Item(4,4).Value = " "
$ws.Cells.Item(4,5).Value = " "
$ws.Cells.Item(4,6).Value = " "
$ws.Cells.Item(4,7).Value = " "
$ws.Cells.Item(4,8).Value = " "
$ws.Cells.Item(4,9).Value = 748804

$ws.Cells.Item(5,4).Value = " "
$ws.Cells.Item(5,5).Value = " "
$ws.Cells.Item(5,6).Value = " "
$ws.Cells.Item(5,7).Value = " "
$ws.Cells.Item(5,8).Value = " "
$ws.Cells.Item(5,9).Value = 861710

$ws.Cells.Item(6,4).Value = " "
$ws.Cells.Item(6,5).Value = " "
$ws.Cells.Item(6,6).Value = " "
$ws.Cells.Item(6,7).Value = " "
$ws.Cells.Item(6,8).Value = " "
$ws.Cells.Item(6,9).Value = 181511

$ws.Cells.Item(7,4).Value = " "
$ws.Cells.Item(7,5).Value = " "
$ws.Cells.Item(7,6).Value = " "
$ws.Cells.Item(7,7).Value = " "
$ws.Cells.Item(7,8).Value = " "
$ws.Cells.Item(7,9).Value = 165856

# --- Append the new client submitted via the "Add New Client" form ---
$newRow = 8
$ws.Cells.Item($newRow,1).Value = "Aaron"
$ws.Cells.Item($newRow,2).Value = "Gallaway"
$ws.Cells.Item($newRow,3).Value = "aaron@gallaway.us"
$ws.Cells.Item($newRow,4).Value = "52 valley street"
$ws.Cells.Item($newRow,5).Value = "providence"
$ws.Cells.Item($newRow,6).Value = "ri"
$ws.Cells.Item($newRow,7).Value = 2909
$ws.Cells.Item($newRow,8).Value = "401-378-6008"
$ws.Cells.Item($newRow,9).Value = 931741

# --- Select the first cell of the newly added row, as the form leaves it ---
$ws.Range("A8").Select()
